# Update cryptos list values per the commit diff (prices + volume % changes,
# plus an OKB/TheGraph row swap at rows 40-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.612.49"
$ws.Range("E2").Value = "  +2.35%  "

# Row 3
$ws.Range("D3").Value = "3.083.38"
$ws.Range("E3").Value = "  +4.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'578.24"
$ws.Range("E5").Value = "  +1.53%  "

# Row 6
$ws.Range("D6").Value = "'168.33"
$ws.Range("E6").Value = "  +5.57%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "3.079.17"
$ws.Range("E8").Value = "  +4.58%  "

# Row 9
$ws.Range("E9").Value = "  +0.79%  "

# Row 10
$ws.Range("D10").Value = "'6.56"
$ws.Range("E10").Value = "  -1.36%  "

# Row 11
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  +1.85%  "

# Row 12
$ws.Range("E12").Value = "  +5.39%  "

# Row 14
$ws.Range("D14").Value = "'36.42"
$ws.Range("E14").Value = "  +6.88%  "

# Row 15
$ws.Range("E15").Value = "  -0.66%  "

# Row 16
$ws.Range("D16").Value = "3.590.63"
$ws.Range("E16").Value = "  +4.28%  "

# Row 17
$ws.Range("D17").Value = "66.650.16"
$ws.Range("E17").Value = "  +2.31%  "

# Row 18
$ws.Range("D18").Value = "'7.21"
$ws.Range("E18").Value = "  +4.17%  "

# Row 19
$ws.Range("D19").Value = "3.080.45"
$ws.Range("E19").Value = "  +4.30%  "

# Row 20
$ws.Range("D20").Value = "'16.18"
$ws.Range("E20").Value = "  +15.05%  "

# Row 21
$ws.Range("D21").Value = "'466.18"
$ws.Range("E21").Value = "  +4.54%  "

# Row 22
$ws.Range("D22").Value = "'0.714"
$ws.Range("E22").Value = "  +4.89%  "

# Row 23
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  +3.93%  "

# Row 24
$ws.Range("D24").Value = "'83.10"
$ws.Range("E24").Value = "  +0.83%  "

# Row 25
$ws.Range("D25").Value = "'2.32"
$ws.Range("E25").Value = "  +4.50%  "

# Row 26
$ws.Range("D26").Value = "'12.84"
$ws.Range("E26").Value = "  +6.54%  "

# Row 27
$ws.Range("D27").Value = "'10.13"
$ws.Range("E27").Value = "  +0.63%  "

# Row 28
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("D29").Value = "'7.99"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$ws.Range("D30").Value = "'2.41"
$ws.Range("E30").Value = "  +0.48%  "

# Row 31
$ws.Range("E31").Value = "  +3.28%  "

# Row 32
$ws.Range("E32").Value = "  -0.21%  "

# Row 33
$ws.Range("D33").Value = "'28.16"
$ws.Range("E33").Value = "  +3.69%  "

# Row 34
$ws.Range("E34").Value = "  +3.95%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.07%  "

# Row 36
$ws.Range("E36").Value = "  +2.81%  "

# Row 37
$ws.Range("D37").Value = "'5.88"
$ws.Range("E37").Value = "  +3.45%  "

# Row 38
$ws.Range("D38").Value = "'2.13"
$ws.Range("E38").Value = "  +7.43%  "

# Row 39
$ws.Range("D39").Value = "'46.72"
$ws.Range("E39").Value = "  +6.46%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'50.26"
$ws.Range("E40").Value = "  +2.66%  "

# Row 41
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.316"
$ws.Range("E41").Value = "  +6.05%  "

# Row 42
$ws.Range("E42").Value = "  +1.98%  "

# Row 43
$ws.Range("D43").Value = "'8.68"
$ws.Range("E43").Value = "  +3.21%  "

# Row 44
$ws.Range("D44").Value = "'2.82"
$ws.Range("E44").Value = "  -0.69%  "

# Row 45
$ws.Range("E45").Value = "  +2.79%  "

# Row 46
$ws.Range("D46").Value = "'381.94"
$ws.Range("E46").Value = "  -1.36%  "

# Row 47
$ws.Range("D47").Value = "2.753.84"
$ws.Range("E47").Value = "  +1.30%  "

# Row 48
$ws.Range("E48").Value = "  +1.71%  "

# Row 49
$ws.Range("E49").Value = "  +0.00%  "

# Row 50
$ws.Range("D50").Value = "'24.59"
$ws.Range("E50").Value = "  +6.04%  "

# Row 51
$ws.Range("D51").Value = "'2.24"
$ws.Range("E51").Value = "  +3.56%  "
